# Prototypes: #4 wget explanation generator.
# Adds three new "wget" tutorial rows (23-25) to the tutorials table, each
# with a "Link" hyperlink in column B, matching data in C:G, and leaves the
# selection on the first empty row below the new data (B26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing 22 data rows (sheet rows 2..22).
$rows = @(
    @{ Id = 22; Topic = "wget"; Name = "All the Wget Commands You Should Know";
       Purpose = "Teach"; Verbosity = "Verbose";
       Glossed = "terms: mirror; number expansion {}, Terminal start, refer, cookies";
       Url = "https://www.maketecheasier.com/all-wget-commands/" },
    @{ Id = 23; Topic = "Terminal scraping"; Name = "Using the Linux Shell for Web Scraping";
       Purpose = "Teach"; Verbosity = "Succinct";
       Glossed = "wget purpose, piping, -O -i options, CSS selector";
       Url = "https://www.linode.com/docs/guides/using-the-linux-shell-for-web-scraping/" },
    @{ Id = 24; Topic = "Terminal scraping"; Name = "Downloading an Entire Web Site with wget";
       Purpose = "Teach"; Verbosity = "Succinct";
       Glossed = "Purpose of wget, where to run command, optimizations";
       Url = "https://www.guyrutenberg.com/2014/05/02/downloading-an-entire-web-site-with-wget/" }
)

$startRow = 23
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.Id
    $ws.Cells.Item($r, 2).Value = "Link"
    $ws.Cells.Item($r, 3).Value = $data.Topic
    $ws.Cells.Item($r, 4).Value = $data.Name
    $ws.Cells.Item($r, 5).Value = $data.Purpose
    $ws.Cells.Item($r, 6).Value = $data.Verbosity
    $ws.Cells.Item($r, 7).Value = $data.Glossed

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $data.Url)
    $ws.Cells.Item($r, 2).Font.Size = 12
}

# Leave the selection on the next empty row, as in the authored edit.
$ws.Range("B26").Select()
